# Regenerate the localization-status report:
#  - Status changes from "Ready for handoff" to "In Translation" for the
#    two tracked files, on the Overview sheet (columns E/F, one per locale)
#    and on each per-locale sheet (zh-cn, de-de -> column C "Status").
#  - The Status column narrows to fit the new (shorter) text on each sheet
#    that shows it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) and de-de (F) status columns, rows 2-3 ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), rows 2-3 ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C), rows 2-3 ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Narrow the Status columns now that the text is shorter ---
$overview.Columns("E").ColumnWidth = 16.333333333333336
$overview.Columns("F").ColumnWidth = 16.333333333333336
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
